$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh prepends two new price records (dated 2022-01-24 /
# serial 44585) to the Albahaca - Vega Modelo de Temuco sheet, pushing the
# previously existing rows 178-187 down to 180-189.
$ws.Rows("178:179").Insert()

# New row 178
$ws.Cells.Item(178, 1).Value = 10
$ws.Cells.Item(178, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(178, 3).Value = "La Araucanía"
$ws.Cells.Item(178, 4).Value = 44585
$ws.Cells.Item(178, 5).Value = 9
$ws.Cells.Item(178, 6).Value = 100112052
$ws.Cells.Item(178, 7).Value = "Albahaca"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 35
$ws.Cells.Item(178, 11).Value = 6000
$ws.Cells.Item(178, 12).Value = 6000
$ws.Cells.Item(178, 13).Value = 6000
$ws.Cells.Item(178, 14).Value = "$/paquete"
$ws.Cells.Item(178, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(178, 16).Value = 6000
$ws.Cells.Item(178, 17).Value = 1
$ws.Cells.Item(178, 18).Value = "Hortaliza"

# New row 179
$ws.Cells.Item(179, 1).Value = 10
$ws.Cells.Item(179, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(179, 3).Value = "La Araucanía"
$ws.Cells.Item(179, 4).Value = 44585
$ws.Cells.Item(179, 5).Value = 9
$ws.Cells.Item(179, 6).Value = 100112052
$ws.Cells.Item(179, 7).Value = "Albahaca"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 55
$ws.Cells.Item(179, 11).Value = 5000
$ws.Cells.Item(179, 12).Value = 5000
$ws.Cells.Item(179, 13).Value = 5000
$ws.Cells.Item(179, 14).Value = "$/paquete"
$ws.Cells.Item(179, 15).Value = "Región del Maule"
$ws.Cells.Item(179, 16).Value = 5000
$ws.Cells.Item(179, 17).Value = 1
$ws.Cells.Item(179, 18).Value = "Hortaliza"
